$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from the last existing data row (A53) onto the
# two new rows so the new date cells render the same way (numFmtId 14).
$ws.Cells.Item(53, 1).Copy()
$ws.Range("A54:A55").PasteSpecial(-4122)

# Row 54: 5/3/2023 - Your First Solo - Thunderstorms and Convective Forecasts
$ws.Cells.Item(54, 1).Value = 45049
$ws.Cells.Item(54, 2).Value = "Your First Solo "

# Row 55: 5/4/2023 - Your First Solo - Radar Imagery
$ws.Cells.Item(55, 1).Value = 45050
$ws.Cells.Item(55, 2).Value = "Your First Solo "

# Write "Radar Imagery" before "Thunderstorms and Convective Forecasts" so the
# shared-string table picks up the same ordering (59 = Radar Imagery,
# 60 = Thunderstorms and Convective Forecasts) as the target workbook.
$ws.Cells.Item(55, 3).Value = "Radar Imagery"
$ws.Cells.Item(54, 3).Value = "Thunderstorms and Convective Forecasts"

# Widen column C slightly (manual resize, no longer auto "best fit").
$ws.Columns.Item(3).ColumnWidth = 36.83

# Update the selection to match where the author left off.
$ws.Range("G48").Select()
